# Apply the "VCU - ajout protocoles de tests / interface microSD" update.
#
# Summary of changes on sheet "Composants utilisés":
#   - New email contact hyperlink in F1 (mathieu.astagneau@etu.ec-lyon.fr)
#   - New token/value in F2 (8vr91L2xhmMpjdtA5ubI)
#   - Two new component rows (7 & 8): microSD connector + logic level shifter,
#     each with hyperlinks on the reference (col B) and the "Datasheet" label (col C)
#   - New column F width
# Sheet "Correspondance PIN VCU" content is untouched (its shared-string
# indices shift automatically because new strings were introduced earlier
# in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1. New contact e-mail in F1 (with mailto hyperlink), and token in F2
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 6).Value2 = "mathieu.astagneau@etu.ec-lyon.fr"
$ws.Hyperlinks.Add($ws.Cells.Item(1, 6), "mailto:mathieu.astagneau@etu.ec-lyon.fr", [System.Type]::Missing, [System.Type]::Missing, "mathieu.astagneau@etu.ec-lyon.fr") | Out-Null

# Match the existing hyperlink look (blue Arial, same as column B/C links)
$ws.Range("B2").Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(2, 6).Value2 = "8vr91L2xhmMpjdtA5ubI"

# New column F width
$ws.Columns.Item(6).ColumnWidth = 37

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Two new component rows
# ---------------------------------------------------------------------
$ws.Cells.Item(7, 1).Value2 = "Connecteur Carte microSD"
$ws.Cells.Item(7, 2).Value2 = "DM3D-SF"
$ws.Cells.Item(7, 3).Value2 = "Datasheet"

$ws.Cells.Item(8, 1).Value2 = "Convertisseur logique 5V-3,3V"
$ws.Cells.Item(8, 2).Value2 = "SN74LV1T34DCKRG4"
$ws.Cells.Item(8, 3).Value2 = "Datasheet"

$ws.Hyperlinks.Add($ws.Cells.Item(7, 2), "https://www.mouser.fr/ProductDetail/Hirose-Connector/DM3D-SF", [System.Type]::Missing, [System.Type]::Missing, "DM3D-SF") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(7, 3), "https://www.hirose.com/product/document?clcode=DM3D-SF&productname=DM3D-SF&series=DM3&documenttype=Catalog&lang=en", [System.Type]::Missing, [System.Type]::Missing, "Datasheet") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(8, 2), "https://www.ti.com/product/SN74LV1T34", [System.Type]::Missing, [System.Type]::Missing, "SN74LV1T34DCKRG4") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(8, 3), "https://www.ti.com/lit/ds/symlink/sn74lv1t34.pdf", [System.Type]::Missing, [System.Type]::Missing, "Datasheet") | Out-Null

# Re-apply the standard reference/datasheet look to the new rows
$ws.Range("B2:C2").Copy()
$ws.Range("B7:C8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. View/selection state (cosmetic, mirrors the authored workbook)
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("E11:E12").Select()

$ws2.Activate()
$ws2.Range("C1").Select()
